$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank rows at the very top of Sheet1, pushing the header
# row (was row 1) and the first data row (was row 2) down to rows 3 and 4.
$ws.Rows("1:2").Insert() | Out-Null

# Re-create the frozen pane so that it now splits after row 3 (the new
# header row) with the top-left visible cell of the scrolling area at A4.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A4").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Leave the selection on the (now blank) first row, selected in full.
$ws.Rows("1:1").Select() | Out-Null

# Reposition the workbook window on screen (reflected as xWindow/yWindow).
$win = $wb.Windows.Item(1)
$win.Left = 3036
$win.Top = 3036
